$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 6;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 18;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 25;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 40;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 54;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 59;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 72;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 85;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 101; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 103; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 110; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 115; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 119; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 122; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.Tag
    $ws.Cells.Item($change.Row, 10).Value = $change.Act
}
